# Add new element-property data (elastic modulus "E" and "E_type") for
# elements 2 and 3 (rows 3 and 4), matching the values already present
# for element 1 (row 2), and update the active cell selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (element 1) already carries the formatting we want for the new
# cells, so copy its number-format / font down into rows 3-4 for the
# columns being populated before writing the new values.
$ws.Range("G2").Copy()
$ws.Range("G3:G4").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("L2").Copy()
$ws.Range("L3:L4").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("G3").Value = 300
$ws.Range("L3").Value = 9

$ws.Range("G4").Value = 300
$ws.Range("L4").Value = 9

$excel.CutCopyMode = $false

# Matches the saved selection state recorded in the sheet (activeCell
# moved from N8 to G5).
$ws.Range("G5").Select()
